$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (per diff) ---
$ws.Range("N13").Value = 46069
$ws.Range("Q13").Value = 212000
$ws.Range("R13").Value = 208000
$ws.Range("S13").Value = 229000
$ws.Range("T13").Value = 232000
$ws.Range("U13").Value = 209000
$ws.Range("N14").Value = 46062
$ws.Range("Q14").Value = 1833000
$ws.Range("R14").Value = 1864000
$ws.Range("S14").Value = 1852000
$ws.Range("T14").Value = 1841000
$ws.Range("U14").Value = 1819000
$ws.Range("F28").Value = -0.01358905968147273
$ws.Range("G28").Value = 0.05441037812530491
$ws.Range("F29").Value = 0.1009743287938078
$ws.Range("G29").Value = 0.1247312200873968
$ws.Range("N29").Value = 46078
$ws.Range("Q29").Value = 2.14
$ws.Range("R29").Value = 2.12
$ws.Range("S29").Value = 2.12
$ws.Range("T29").Value = 2.13
$ws.Range("U29").Value = 2.15
$ws.Range("F30").Value = -0.02377629691120575
$ws.Range("G30").Value = 0.06592882125886601
$ws.Range("N30").Value = 46078
$ws.Range("Q30").Value = 2.28
$ws.Range("R30").Value = 2.26
$ws.Range("S30").Value = 2.26
$ws.Range("T30").Value = 2.28
$ws.Range("U30").Value = 2.29
$ws.Range("F31").Value = 0.09596509959750062
$ws.Range("G31").Value = 0.1265306873826442
$ws.Range("C46").Value = 45992
$ws.Range("C47").Value = 45992
$ws.Range("N47").Value = 46077
$ws.Range("C48").Value = 45992
$ws.Range("N48").Value = 46077
$ws.Range("Q48").Value = 3.43
$ws.Range("R48").Value = 3.43
$ws.Range("S48").Value = 3.48
$ws.Range("T48").Value = 3.47
$ws.Range("U48").Value = 3.47
$ws.Range("C49").Value = 45992
$ws.Range("N49").Value = 46077
$ws.Range("Q49").Value = 3.61
$ws.Range("R49").Value = 3.59
$ws.Range("S49").Value = 3.65
$ws.Range("T49").Value = 3.65
$ws.Range("U49").Value = 3.66
$ws.Range("C50").Value = 45992
$ws.Range("N50").Value = 46077
$ws.Range("Q50").Value = 4.04
$ws.Range("R50").Value = 4.03
$ws.Range("S50").Value = 4.08
$ws.Range("T50").Value = 4.08
$ws.Range("U50").Value = 4.09
$ws.Range("C51").Value = 45992
$ws.Range("N51").Value = 46069
$ws.Range("N52").Value = 46077
$ws.Range("Q52").Value = 5.77
$ws.Range("R52").Value = 5.76
$ws.Range("S52").Value = 5.77
$ws.Range("T52").Value = 5.76
$ws.Range("U52").Value = 5.76

# --- Update formatting: remove yellow highlight fill (style 48 -> 49) ---
# Source cell C28 already carries the target "no fill" direct format (style 49).
$ws.Range("C28").Copy()
$ws.Range("C46").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C47").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C48").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C49").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C50").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C51").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N51").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
